# Doing Updates for Financials
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Income Statement
$ws.Range("D8").Value = 3500    # Total Revenue
$ws.Range("D9").Value = 3200    # Cost of Revenue
$ws.Range("D17").Value = 5600   # Total Operating Expenses
$ws.Range("D18").Value = -2100  # Operating Income or Loss
$ws.Range("D21").Value = -2000  # Earnings Before Interest And Taxes
$ws.Range("D23").Value = -2300  # Income Before Tax
$ws.Range("D26").Value = -2300  # Income After Tax
$ws.Range("D27").Value = -2300  # Net Income From Continuing Ops
$ws.Range("D33").Value = -2300  # Net Income
$ws.Range("D35").Value = -2300  # Net Income Applicable To Common Shares

# Balance Sheet
$ws.Range("D49").Value = 1600   # Goodwill
$ws.Range("D54").Value = 2800   # Total Assets
$ws.Range("D60").Value = 4800   # Total Current Liabilities
$ws.Range("D66").Value = 6300   # Total Liabilities

# Cash Flow Statement
$ws.Range("D72").Value = -9200  # Retained Earnings
$ws.Range("D76").Value = -3500  # Total Stockholder Equity
$ws.Range("D81").Value = -2300  # Net Income
$ws.Range("D89").Value = -1500  # Total Cash Flow From Operating Activities
$ws.Range("H91").Value = "NA"   # Capital Expenditures
